# Append two new match rows (5 and 6) to the "Jogos da Semana" sheet,
# matching the FlashScore export format already used by rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row5Data = @("nVFWn41R", "09/10/2024", "23:07", "USA - USL CHAMPIONSHIP", "Orange County SC", "Colorado Springs", 2.45, 3.25, 2.75, 3.1, 2.07, 3.3, 1.06, 7.2, 1.31, 3.2, 1.93, 1.82, 1.42, 2.65, 1.7, 2.02, 8.25, 12, 9.25, 26, 20, 29, 7.2, 6.2, 13, 60, 450, 9.25, 14.5, 10, 32, 22, 29, 4.4, 13.5, 22, 60, 100, 300, 2.65, 7, 65, 4.7, 15, 22, 65, 100, 300, 51, 51)
$row6Data = @("lv2brO9r", "09/10/2024", "23:11", "USA - USL CHAMPIONSHIP", "Sacramento Republic", "Monterey Bay", 1.44, 4.2, 6.8, 1.98, 2.25, 6.4, 1.05, 7.6, 1.28, 3.35, 1.85, 1.88, 1.4, 2.75, 2.05, 1.7, 6.2, 6.3, 8.25, 9.25, 12.5, 30, 7.6, 8.25, 21, 110, 1000, 16, 40, 21, 150, 80, 75, 3.15, 6.8, 18.5, 21, 55, 300, 2.75, 8.75, 90, 7.8, 40, 45, 300, 350, 450, 51, 51)

$rows = @(5, 6)
$data = @($row5Data, $row6Data)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $rows[$r]
    $values = $data[$r]
    for ($c = 1; $c -le $values.Length; $c++) {
        $cell = $ws.Cells.Item($rowNum, $c)
        # Column B holds dd/mm/yyyy-looking text ("Date"); force text format
        # first so Excel doesn't silently reinterpret it as a date serial,
        # matching how the source rows store it as plain text.
        if ($c -eq 2) {
            $cell.NumberFormat = "@"
            $cell.Value = $values[$c - 1]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $values[$c - 1]
        }
    }
}
